{"js": "// [Todo] Increase priority of item.\n//\n// The Todo list uses font color as an informal priority marker\n// (red = FF0000 for high priority, green = 00B050 for lower priority).\n// The item \"Improve genericity of parameter passing in remote function\n// caller.\" is being bumped to high priority, so its text (and the\n// paragraph mark) switch from green to red.\n\nconst body = context.document.body;\n\n// \"of parameter passing in remote function caller\" is unique in the\n// document (the similarly-worded \"Improve API genericity.\" item earlier\n// in the doc does not contain this phrase), so it unambiguously pins the\n// paragraph we need regardless of how its runs are currently split.\nconst results = body.search(\"of parameter passing in remote function caller\", {\n  matchCase: false,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly 1 match for the target Todo item, found \" + results.items.length\n  );\n}\n\nconst paragraph = results.items[0].paragraphs.getFirst();\nparagraph.load(\"text\");\nawait context.sync();\n\nif (paragraph.text.indexOf(\"genericity\") === -1) {\n  throw new Error(\"Unexpected paragraph located: \" + paragraph.text);\n}\n\n// Setting color on the paragraph (rather than just the matched run) also\n// recolors the paragraph mark's rPr and every run in the paragraph\n// (\"Improve \", \"genericity\", \" of parameter passing in remote function\n// caller.\"), matching the diff exactly.\nparagraph.font.color = \"#FF0000\";\n\nawait context.sync();\n", "ps1": "# [Todo] Increase priority of item.\n#\n# The Todo list uses font color as an informal priority marker\n# (red = FF0000 for high priority, green = 00B050 for lower priority).\n# The item \"Improve genericity of parameter passing in remote function\n# caller.\" is being bumped to high priority, so its text (and the\n# paragraph mark) switch from green to red.\n\n$d = $word.ActiveDocument\n\n# \"of parameter passing in remote function caller\" is unique in the\n# document (the similarly-worded \"Improve API genericity.\" item earlier\n# in the doc does not contain this phrase), so it unambiguously pins the\n# paragraph we need regardless of how its runs are currently split.\n$needle = \"of parameter passing in remote function caller\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$needle*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the target Todo item paragraph.\"\n}\n\n# Paint the whole paragraph (all of its runs, plus the paragraph mark's\n# own run properties) red -- wdColorRed (255 == RGB(255,0,0) == 0xFF0000).\n$target.Range.Font.Color = 255\n"}
